$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update confusion matrix values (corrected translations)
$ws.Range("C3").Value = "20 (0.9624)"
$ws.Range("D3").Value = "1 (0.0476)"
$ws.Range("C4").Value = "10 (0.7143)"
$ws.Range("D4").Value = "4 (0.2857)"

# Column C needs to widen to fit the new, longer text
$ws.Columns.Item(3).ColumnWidth = 9.7

# Update the selected cell shown when the sheet was last saved
$ws.Range("F4").Select() | Out-Null
